# Auto-generated edit script: updates Leve price/profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{cell="H9"; value=716.7},
    @{cell="I9"; value=782.4286},
    @{cell="K9"; value=782.4286},
    @{cell="M9"; value=-613.4286},
    @{cell="H32"; value=10898},
    @{cell="I32"; value=13456.333},
    @{cell="J32"; value=8979.25},
    @{cell="K32"; value=13456.333},
    @{cell="L32"; value=8979.25},
    @{cell="M32"; value=-13130.333},
    @{cell="N32"; value=-9631.25},
    @{cell="H74"; value=5963.3335},
    @{cell="I74"; value=4924.2856},
    @{cell="K74"; value=4924.2856},
    @{cell="M74"; value=-3988.2856},
    @{cell="H77"; value=5963.3335},
    @{cell="I77"; value=4924.2856},
    @{cell="K77"; value=24621.428},
    @{cell="M77"; value=-19941.428},
    @{cell="H94"; value=1000},
    @{cell="I94"; value=1000},
    @{cell="K94"; value=1000},
    @{cell="M94"; value=-549},
    @{cell="H95"; value=48656},
    @{cell="J95"; value=48656},
    @{cell="L95"; value=48656},
    @{cell="N95"; value=-54148},
    @{cell="H100"; value=2162.2},
    @{cell="I100"; value=2315.25},
    @{cell="J100"; value=1550},
    @{cell="K100"; value=2315.25},
    @{cell="L100"; value=1550},
    @{cell="M100"; value=-1774.25},
    @{cell="N100"; value=-2632},
    @{cell="H106"; value=12299.091},
    @{cell="I106"; value=13255.223},
    @{cell="J106"; value=7996.5},
    @{cell="K106"; value=13255.223},
    @{cell="L106"; value=7996.5},
    @{cell="M106"; value=-12624.223},
    @{cell="N106"; value=-9258.5},
    @{cell="H111"; value=2599.8333},
    @{cell="I111"; value=2599.8333},
    @{cell="J111"; value=0},
    @{cell="K111"; value=7799.499899999999},
    @{cell="L111"; value=0},
    @{cell="M111"; value=-4732.499899999999},
    @{cell="N111"; value=$null},
    @{cell="H113"; value=3018.6155},
    @{cell="I113"; value=2324.3},
    @{cell="K113"; value=2324.3},
    @{cell="M113"; value=929.6999999999998},
    @{cell="H116"; value=5330.3335},
    @{cell="I116"; value=5330.3335},
    @{cell="K116"; value=5330.3335},
    @{cell="M116"; value=-1888.3335},
    @{cell="H129"; value=2320.125},
    @{cell="I129"; value=1495.5555},
    @{cell="K129"; value=4486.666499999999},
    @{cell="M129"; value=513.3335000000006},
    @{cell="H135"; value=0},
    @{cell="I135"; value=0},
    @{cell="K135"; value=0},
    @{cell="M135"; value=$null}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{cell="H2"; value=888},
    @{cell="I2"; value=888},
    @{cell="K2"; value=888},
    @{cell="M2"; value=-775},
    @{cell="H5"; value=287.42856},
    @{cell="J5"; value=198},
    @{cell="L5"; value=198},
    @{cell="N5"; value=-422},
    @{cell="H14"; value=14750},
    @{cell="J14"; value=15000},
    @{cell="L14"; value=15000},
    @{cell="N14"; value=-15350},
    @{cell="H32"; value=6842.0713},
    @{cell="I32"; value=6022.269},
    @{cell="K32"; value=6022.269},
    @{cell="M32"; value=-5735.269},
    @{cell="H45"; value=2359.2},
    @{cell="I45"; value=2199.25},
    @{cell="K45"; value=2199.25},
    @{cell="M45"; value=-1822.25},
    @{cell="H63"; value=5828.909},
    @{cell="I63"; value=4947},
    @{cell="K63"; value=4947},
    @{cell="M63"; value=-4261},
    @{cell="H66"; value=5828.909},
    @{cell="I66"; value=4947},
    @{cell="K66"; value=24735},
    @{cell="M66"; value=-21303},
    @{cell="H97"; value=1281.0646},
    @{cell="I97"; value=675.5},
    @{cell="K97"; value=675.5},
    @{cell="M97"; value=-179.5},
    @{cell="H102"; value=2286.875},
    @{cell="I102"; value=2040.5714},
    @{cell="K102"; value=2040.5714},
    @{cell="M102"; value=-418.5714},
    @{cell="H110"; value=6345.9565},
    @{cell="I110"; value=7946.1665},
    @{cell="J110"; value=4600.273},
    @{cell="K110"; value=7946.1665},
    @{cell="L110"; value=4600.273},
    @{cell="M110"; value=-5901.1665},
    @{cell="N110"; value=-8690.273000000001},
    @{cell="H116"; value=888},
    @{cell="I116"; value=888},
    @{cell="K116"; value=888},
    @{cell="M116"; value=1406}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{cell="H3"; value=888},
    @{cell="I3"; value=888},
    @{cell="K3"; value=888},
    @{cell="M3"; value=-774},
    @{cell="H4"; value=287.42856},
    @{cell="J4"; value=198},
    @{cell="L4"; value=198},
    @{cell="N4"; value=-428},
    @{cell="H22"; value=0},
    @{cell="I22"; value=0},
    @{cell="K22"; value=0},
    @{cell="M22"; value=$null},
    @{cell="H99"; value=3609.5},
    @{cell="I99"; value=3978},
    @{cell="K99"; value=3978},
    @{cell="M99"; value=-2480},
    @{cell="H105"; value=2861.5},
    @{cell="I105"; value=2798.75},
    @{cell="J105"; value=2924.25},
    @{cell="K105"; value=2798.75},
    @{cell="L105"; value=2924.25},
    @{cell="M105"; value=-1051.75},
    @{cell="N105"; value=-6418.25},
    @{cell="H107"; value=2302},
    @{cell="I107"; value=2302},
    @{cell="K107"; value=2302},
    @{cell="M107"; value=-382}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{cell="H7"; value=900},
    @{cell="I7"; value=800},
    @{cell="K7"; value=800},
    @{cell="M7"; value=-687},
    @{cell="H22"; value=2858663},
    @{cell="J22"; value=5715663},
    @{cell="L22"; value=5715663},
    @{cell="N22"; value=-5716363},
    @{cell="H122"; value=2418.8333},
    @{cell="I122"; value=2262.6},
    @{cell="K122"; value=6787.799999999999},
    @{cell="M122"; value=-4337.799999999999}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{cell="H122"; value=15969.286},
    @{cell="I122"; value=50997.5},
    @{cell="J122"; value=1958},
    @{cell="K122"; value=458977.5},
    @{cell="L122"; value=17622},
    @{cell="M122"; value=-456527.5},
    @{cell="N122"; value=-22522},
    @{cell="H125"; value=29998},
    @{cell="I125"; value=29998},
    @{cell="J125"; value=0},
    @{cell="K125"; value=89994},
    @{cell="L125"; value=0},
    @{cell="M125"; value=-85074},
    @{cell="N125"; value=$null}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{cell="H70"; value=7488.125},
    @{cell="J70"; value=10500},
    @{cell="L70"; value=10500},
    @{cell="N70"; value=-11040},
    @{cell="H73"; value=7488.125},
    @{cell="J73"; value=10500},
    @{cell="L73"; value=10500},
    @{cell="N73"; value=-12372},
    @{cell="H80"; value=5000},
    @{cell="I80"; value=0},
    @{cell="K80"; value=0},
    @{cell="M80"; value=$null},
    @{cell="H83"; value=5000},
    @{cell="I83"; value=0},
    @{cell="K83"; value=0},
    @{cell="M83"; value=$null},
    @{cell="H107"; value=2300},
    @{cell="I107"; value=375.25},
    @{cell="K107"; value=375.25},
    @{cell="M107"; value=1544.75},
    @{cell="H113"; value=0},
    @{cell="I113"; value=0},
    @{cell="K113"; value=0},
    @{cell="M113"; value=$null},
    @{cell="H136"; value=53796.777},
    @{cell="J136"; value=53796.777},
    @{cell="L136"; value=161390.331},
    @{cell="N136"; value=-166490.331}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{cell="H62"; value=0},
    @{cell="I62"; value=0},
    @{cell="K62"; value=0},
    @{cell="M62"; value=$null},
    @{cell="H65"; value=0},
    @{cell="I65"; value=0},
    @{cell="K65"; value=0},
    @{cell="M65"; value=$null},
    @{cell="H70"; value=34000},
    @{cell="J70"; value=34000},
    @{cell="L70"; value=34000},
    @{cell="N70"; value=-34540},
    @{cell="H73"; value=34000},
    @{cell="J73"; value=34000},
    @{cell="L73"; value=34000},
    @{cell="N73"; value=-35872},
    @{cell="H136"; value=7184.75},
    @{cell="I136"; value=4298.8},
    @{cell="J136"; value=11994.667},
    @{cell="K136"; value=12896.4},
    @{cell="L136"; value=35984.001},
    @{cell="M136"; value=-10346.4},
    @{cell="N136"; value=-41084.001}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{cell="H62"; value=6632.9},
    @{cell="J62"; value=7328.625},
    @{cell="L62"; value=7328.625},
    @{cell="N62"; value=-8576.625},
    @{cell="H65"; value=6632.9},
    @{cell="J65"; value=7328.625},
    @{cell="L65"; value=36643.125},
    @{cell="N65"; value=-42883.125},
    @{cell="H96"; value=6183.4},
    @{cell="J96"; value=6062.75},
    @{cell="L96"; value=6062.75},
    @{cell="N96"; value=-8808.75},
    @{cell="H107"; value=975},
    @{cell="I107"; value=975},
    @{cell="J107"; value=0},
    @{cell="K107"; value=2925},
    @{cell="L107"; value=0},
    @{cell="M107"; value=-1005},
    @{cell="N107"; value=$null},
    @{cell="H126"; value=3787.5264},
    @{cell="I126"; value=3574.4},
    @{cell="J126"; value=4586.75},
    @{cell="K126"; value=10723.2},
    @{cell="L126"; value=13760.25},
    @{cell="M126"; value=-8253.200000000001},
    @{cell="N126"; value=-18700.25}
)
foreach ($u in $updates) {
    if ($u.value -eq $null) {
        $ws.Range($u.cell).ClearContents()
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

